$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regenerated s_val data (filtered save games) for rows 2-9, columns B-E and G.
$data = @{
    2 = @{ B = 1.445647641019636;  C = 1.626987699542094;  D = 0.1496068669990043; E = 0.5333859586016987; G = 3.755628166162433 }
    3 = @{ B = 0.6545652718822623; C = 0.3048912486333797;  D = 0.7210945179870265; E = 0.5333859586016987; G = 2.213936997104367 }
    4 = @{ B = 1.445647641019636;  C = 1.626987699542094;  D = 0.7210945179870265; E = 0.5333859586016987; G = 4.327115817150455 }
    5 = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 0.1496068669990043; E = 0.5333859586016987; G = 5.582307763322248 }
    6 = @{ B = 0.1169995834814548; C = 0.3048912486333797;  D = 0.7210945179870265; E = 0.5333859586016987; G = 1.67637130870356 }
    7 = @{ B = 1.445647641019636;  C = 1.626987699542094;  D = 0.7210945179870265; E = 0.5333859586016987; G = 4.327115817150455 }
    8 = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 0.7210945179870265; E = 0.5333859586016987; G = 6.15379541431027 }
    9 = @{ B = 1.445647641019636;  C = 1.626987699542094;  D = 3.223369029078222;  E = 0.5333859586016987; G = 6.82939032824165 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
